$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.006.88'
$ws.Range('E2').Value = '  +0.55%  '

$ws.Range('D3').Value = '1.911.50'
$ws.Range('E3').Value = '  +0.90%  '

$ws.Range('D4').Value = '0.9977'
$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').Value = '0.8071'
$ws.Range('E5').Value = '  +6.27%  '

$ws.Range('D6').Value = '242.12'
$ws.Range('E6').Value = '  +1.10%  '

$ws.Range('D7').Value = '0.9990'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = '0.3168'
$ws.Range('E8').Value = '  +4.02%  '

$ws.Range('D9').Value = '26.50'
$ws.Range('E9').Value = '  +4.57%  '

$ws.Range('D10').Value = '0.06924'
$ws.Range('E10').Value = '  +1.50%  '

$ws.Range('D11').Value = '0.07997'
$ws.Range('E11').Value = '  +0.06%  '

$ws.Range('D12').Value = '1.910.43'
$ws.Range('E12').Value = '  +1.39%  '

$ws.Range('D13').Value = '0.7413'
$ws.Range('E13').Value = '  -1.05%  '

$ws.Range('D14').Value = '5.202'
$ws.Range('E14').Value = '  -0.02%  '

$ws.Range('D15').Value = '93.08'
$ws.Range('E15').Value = '  +2.20%  '

$ws.Range('D16').Value = '29.989.25'
$ws.Range('E16').Value = '  +0.49%  '

$ws.Range('D17').Value = '14.03'
$ws.Range('E17').Value = '  +1.08%  '

$ws.Range('D18').Value = '5.890'
$ws.Range('E18').Value = '  -1.76%  '

$ws.Range('D19').Value = '246.31'
$ws.Range('E19').Value = '  +4.76%  '

$ws.Range('D20').Value = '0.000007751'
$ws.Range('E20').Value = '  +1.05%  '

$ws.Range('D21').Value = '0.9995'
$ws.Range('E21').Value = '  -0.03%  '

$ws.Range('D22').Value = '2.150.08'
$ws.Range('E22').Value = '  +0.85%  '

$ws.Range('D23').Value = '0.9974'
$ws.Range('E23').Value = '  -0.23%  '

$ws.Range('D24').Value = '6.844'
$ws.Range('E24').Value = '  -1.46%  '

$ws.Range('D25').Value = '167.98'
$ws.Range('E25').Value = '  +1.69%  '

$ws.Range('D26').Value = '9.221'
$ws.Range('E26').Value = '  -0.27%  '

$ws.Range('D27').Value = '0.1424'
$ws.Range('E27').Value = '  +11.01%  '

$ws.Range('D28').Value = '18.94'
$ws.Range('E28').Value = '  +1.28%  '

$ws.Range('D29').Value = '2.040'
$ws.Range('E29').Value = '  -0.24%  '

$ws.Range('E30').Value = '  +1.60%  '

$ws.Range('D31').Value = '1.517'
$ws.Range('E31').Value = '  +0.30%  '

$ws.Range('D32').Value = '4.316'
$ws.Range('E32').Value = '  +0.81%  '

$ws.Range('D33').Value = '4.091'
$ws.Range('E33').Value = '  +1.71%  '

$ws.Range('D34').Value = '0.05493'
$ws.Range('E34').Value = '  +2.39%  '

$ws.Range('D35').Value = '1.270'
$ws.Range('E35').Value = '  +1.58%  '

$ws.Range('D36').Value = '0.7316'

$ws.Range('D37').Value = '2.718'
$ws.Range('E37').Value = '  +0.22%  '

$ws.Range('D38').Value = '0.01927'
$ws.Range('E38').Value = '  +0.10%  '

$ws.Range('D39').Value = '2.782'
$ws.Range('E39').Value = '  +0.45%  '

$ws.Range('D40').Value = '6.162'
$ws.Range('E40').Value = '  -0.50%  '

$ws.Range('D41').Value = '0.4430'
$ws.Range('E41').Value = '  +0.48%  '

$ws.Range('D42').Value = '72.62'
$ws.Range('E42').Value = '  +0.45%  '

$ws.Range('D43').Value = '0.9993'
$ws.Range('E43').Value = '  -0.06%  '

$ws.Range('D44').Value = '0.8368'
$ws.Range('E44').Value = '  +1.51%  '

$ws.Range('D45').Value = '1.879'
$ws.Range('E45').Value = '  -2.00%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.569'
$ws.Range('E46').Value = '  -0.29%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '100.48'
$ws.Range('E47').Value = '  -0.52%  '

$ws.Range('D48').Value = '9.742'
$ws.Range('E48').Value = '  -0.28%  '

$ws.Range('D49').Value = '985.18'
$ws.Range('E49').Value = '  +6.97%  '

$ws.Range('D50').Value = '2.057.04'
$ws.Range('E50').Value = '  +0.87%  '

$ws.Range('D51').Value = '36.28'
$ws.Range('E51').Value = '  +0.36%  '
